# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.455362044514542;    C = 1.655778082260271;   D = 0.7527432677738641;  E = 10.19245300693656;   G = 14.05633640148523 }
    3 = @{ B = 0.6606524410359556;   C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697;  G = 2.960089034096801 }
    4 = @{ B = 0.6606524410359556;   C = 250555.8564151394;   D = 0.7527432677738641;  E = 1133.036916526867;   G = 251690.306727375 }
    5 = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.1494219747398047;  E = 0.4942365360607697;  G = 5.586269137925634 }
    6 = @{ B = 0.00001292064567892659; C = 0.306821227259698; D = 261.3203778131603;   E = 1133.036916526867;   G = 1394.664128487933 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $vals.G
}
